# Insert a new row at position 110, shifting existing rows 110-154 down to 111-155.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record.
$ws.Cells.Item(110, 1).Value = 6
$ws.Cells.Item(110, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(110, 3).Value = "Metropolitana"
$ws.Cells.Item(110, 4).Value = 44572
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 100112001
$ws.Cells.Item(110, 7).Value = "Berenjena"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 200
$ws.Cells.Item(110, 11).Value = 10000
$ws.Cells.Item(110, 12).Value = 12000
$ws.Cells.Item(110, 13).Value = 11200
$ws.Cells.Item(110, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 187
$ws.Cells.Item(110, 17).Value = 60
$ws.Cells.Item(110, 18).Value = "Hortaliza"
